$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2 and 3 with new capital structure values
foreach ($row in 2..3) {
    $ws.Range("D$row").Value = 0.0495
    $ws.Range("E$row").Value = 0.0444
    $ws.Range("F$row").ClearContents()

    $ws.Range("K$row").Value = 72.2
    $ws.Range("L$row").Value = 0.5197984161267099
    $ws.Range("M$row").Value = 30.5803
    $ws.Range("N$row").Value = 0.07202143193593971
    $ws.Range("O$row").Value = 0.4235498614958449
    $ws.Range("P$row").Value = 15.7803
    $ws.Range("Q$row").Value = 0.03716509656146962
    $ws.Range("R$row").Value = 0.2185637119113573
    $ws.Range("S$row").Value = 14.8
    $ws.Range("T$row").Value = 0.483971707275599
    $ws.Range("U$row").Value = 207.1
    $ws.Range("V$row").Value = 0.487753179463024
    $ws.Range("W$row").Value = 0.0991349718522587
    $ws.Range("X$row").Value = 0.06316325696491433
    $ws.Range("Y$row").Value = 0.03597171488734437
    $ws.Range("Z$row").Value = 0.1380852967491799
    $ws.Range("AB$row").Value = 0.04067772647105903
    $ws.Range("AC$row").Value = -0.04067772647105903
    $ws.Range("AD$row").Value = 415.3
    $ws.Range("AF$row").Value = 415.3
    $ws.Range("AG$row").Value = 208.2
    $ws.Range("AH$row").Value = 0.4944636266222169
    $ws.Range("AI$row").Value = 0.3479681608713867
    $ws.Range("AJ$row").Value = 0.3290139064475348
    $ws.Range("AK$row").Value = 0.2110705596107056
    $ws.Range("AM$row").Value = -0.209
}
